$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 931.2963
$ws.Range("I33").Value = 797.2381
$ws.Range("J33").Value = 1400.5
$ws.Range("K33").Value = 797.2381
$ws.Range("L33").Value = 1400.5
$ws.Range("M33").Value = -568.2381
$ws.Range("N33").Value = -1858.5
$ws.Range("H76").Value = 3999.3333
$ws.Range("I76").Value = 4249.5
$ws.Range("J76").Value = 3499
$ws.Range("K76").Value = 4249.5
$ws.Range("L76").Value = 3499
$ws.Range("M76").Value = -3934.5
$ws.Range("N76").Value = -4129
$ws.Range("H79").Value = 3999.3333
$ws.Range("I79").Value = 4249.5
$ws.Range("J79").Value = 3499
$ws.Range("K79").Value = 4249.5
$ws.Range("L79").Value = 3499
$ws.Range("M79").Value = -3157.5
$ws.Range("N79").Value = -5683
$ws.Range("H98").Value = 47620430
$ws.Range("I98").Value = 52632896
$ws.Range("J98").Value = 1999.5
$ws.Range("K98").Value = 52632896
$ws.Range("L98").Value = 1999.5
$ws.Range("M98").Value = -52631398
$ws.Range("N98").Value = -4995.5
$ws.Range("H113").Value = 45460544
$ws.Range("J113").Value = 80006000
$ws.Range("L113").Value = 80006000
$ws.Range("N113").Value = -80012508
$ws.Range("H122").Value = 47620430
$ws.Range("I122").Value = 52632896
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 157898688
$ws.Range("L122").Value = 5998.5
$ws.Range("M122").Value = -157896238
$ws.Range("N122").Value = -10898.5
$ws.Range("H138").Value = 3219.1177
$ws.Range("J138").Value = 3421.8667
$ws.Range("L138").Value = 10265.6001
$ws.Range("N138").Value = -20545.6001
$ws.Range("H141").Value = 5000
$ws.Range("I141").Value = 5000
$ws.Range("K141").Value = 15000
$ws.Range("M141").Value = -9820

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1228.5454
$ws.Range("I2").Value = 1228.5454
$ws.Range("K2").Value = 1228.5454
$ws.Range("M2").Value = -1115.5454
$ws.Range("H12").Value = 2212.5
$ws.Range("I12").Value = 1200
$ws.Range("J12").Value = 3900
$ws.Range("K12").Value = 1200
$ws.Range("L12").Value = 3900
$ws.Range("M12").Value = -1027
$ws.Range("N12").Value = -4246
$ws.Range("H31").Value = 48622.285
$ws.Range("I31").Value = 2249.5
$ws.Range("J31").Value = 110452.664
$ws.Range("K31").Value = 2249.5
$ws.Range("L31").Value = 110452.664
$ws.Range("M31").Value = -1955.5
$ws.Range("N31").Value = -111040.664
$ws.Range("H32").Value = 11368538
$ws.Range("I32").Value = 11909891
$ws.Range("J32").Value = 122.5
$ws.Range("K32").Value = 11909891
$ws.Range("L32").Value = 122.5
$ws.Range("M32").Value = -11909604
$ws.Range("N32").Value = -696.5
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""
$ws.Range("H116").Value = 1228.5454
$ws.Range("I116").Value = 1228.5454
$ws.Range("K116").Value = 1228.5454
$ws.Range("M116").Value = 1065.4546
$ws.Range("H132").Value = 4887.8125
$ws.Range("I132").Value = 2717.1155
$ws.Range("K132").Value = 8151.3465
$ws.Range("M132").Value = -5621.3465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1228.5454
$ws.Range("I3").Value = 1228.5454
$ws.Range("K3").Value = 1228.5454
$ws.Range("M3").Value = -1114.5454
$ws.Range("H22").Value = 762.75
$ws.Range("I22").Value = 762.75
$ws.Range("K22").Value = 762.75
$ws.Range("M22").Value = -589.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2658.3125
$ws.Range("J58").Value = 2436.5
$ws.Range("L58").Value = 2436.5
$ws.Range("N58").Value = -2842.5
$ws.Range("H62").Value = 2698.5
$ws.Range("I62").Value = 2698
$ws.Range("K62").Value = 2698
$ws.Range("M62").Value = -2074
$ws.Range("H65").Value = 2698.5
$ws.Range("I65").Value = 2698
$ws.Range("K65").Value = 13490
$ws.Range("M65").Value = -10370
$ws.Range("H132").Value = 2000
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 6000
$ws.Range("N132").Value = -11060
$ws.Range("H134").Value = 717562.4
$ws.Range("I134").Value = 1251486.1
$ws.Range("J134").Value = 5664
$ws.Range("K134").Value = 3754458.3
$ws.Range("L134").Value = 16992
$ws.Range("M134").Value = -3751923.3
$ws.Range("N134").Value = -22062
$ws.Range("H136").Value = 2658.3125
$ws.Range("J136").Value = 2436.5
$ws.Range("L136").Value = 7309.5
$ws.Range("N136").Value = -12409.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 424.125
$ws.Range("I2").Value = 542.5
$ws.Range("K2").Value = 3255
$ws.Range("M2").Value = -3142
$ws.Range("H20").Value = 1125
$ws.Range("I20").Value = 500
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = -1273
$ws.Range("N20").Value = -9454
$ws.Range("H141").Value = 10334.059
$ws.Range("I141").Value = 9067.9
$ws.Range("J141").Value = 12142.857
$ws.Range("K141").Value = 27203.7
$ws.Range("L141").Value = 36428.571
$ws.Range("M141").Value = -22023.7
$ws.Range("N141").Value = -46788.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 62508820
$ws.Range("I132").Value = 76932780
$ws.Range("K132").Value = 230798340
$ws.Range("M132").Value = -230795810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3634
$ws.Range("I22").Value = 3634
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3634
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3339
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 3634
$ws.Range("I27").Value = 3634
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3634
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3527
$ws.Range("N27").Value = ""
$ws.Range("H46").Value = 2455.1738
$ws.Range("I46").Value = 2093.1052
$ws.Range("J46").Value = 4175
$ws.Range("K46").Value = 2093.1052
$ws.Range("L46").Value = 4175
$ws.Range("M46").Value = -1905.1052
$ws.Range("N46").Value = -4551
$ws.Range("H68").Value = 4033
$ws.Range("J68").Value = 5967.3335
$ws.Range("L68").Value = 5967.3335
$ws.Range("N68").Value = -7465.3335
$ws.Range("H71").Value = 4033
$ws.Range("J71").Value = 5967.3335
$ws.Range("L71").Value = 29836.6675
$ws.Range("N71").Value = -37324.6675
$ws.Range("H127").Value = 75699.75
$ws.Range("J127").Value = 75699.75
$ws.Range("L127").Value = 75699.75
$ws.Range("N127").Value = -85619.75
$ws.Range("H132").Value = 1015542.75
$ws.Range("I132").Value = 1252310
$ws.Range("K132").Value = 3756930
$ws.Range("M132").Value = -3754400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 56004.668
$ws.Range("J15").Value = 56004.668
$ws.Range("L15").Value = 56004.668
$ws.Range("N15").Value = -56580.668
